$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 4.3
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 1.8
$ws.Range("I2").Value = 1.92
$ws.Range("O2").Value = 1.24
$ws.Range("Q2").Value = 1.73
$ws.Range("T2").Value = 1.76
$ws.Range("V2").Value = 2.08
$ws.Range("W2").Value = 1.25
$ws.Range("X2").Value = 90
$ws.Range("Y2").Value = 16.5
$ws.Range("Z2").Value = 24
$ws.Range("AB2").Value = 970
$ws.Range("AC2").Value = 9.800000000000001
$ws.Range("AH2").Value = 40
$ws.Range("AO2").Value = 29

# Row 3
$ws.Range("L3").Value = 1.01
$ws.Range("AI3").Value = 500
$ws.Range("AL3").Value = 500
$ws.Range("AM3").Value = 500

# Row 4
$ws.Range("F4").Value = 1.75
$ws.Range("G4").Value = 1.87
$ws.Range("H4").Value = 4.9
$ws.Range("I4").Value = 5.7
$ws.Range("P4").Value = 1.98
$ws.Range("V4").Value = 1.21
$ws.Range("W4").Value = 2.14
$ws.Range("X4").Value = 16.5
$ws.Range("Y4").Value = 19.5
$ws.Range("Z4").Value = 44
$ws.Range("AB4").Value = 9.199999999999999
$ws.Range("AD4").Value = 22
$ws.Range("AJ4").Value = 20
$ws.Range("AL4").Value = 500
$ws.Range("AN4").Value = 12

# Row 5
$ws.Range("F5").Value = 1.72
$ws.Range("I5").Value = 6.4
$ws.Range("K5").Value = 3.85
$ws.Range("L5").Value = 1.4
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 3.35
$ws.Range("P5").Value = 1.8
$ws.Range("Q5").Value = 2.08
$ws.Range("S5").Value = 3.75
$ws.Range("V5").Value = 1.18
$ws.Range("X5").Value = 90
$ws.Range("Y5").Value = 500
$ws.Range("AA5").Value = 700
$ws.Range("AC5").Value = 13
$ws.Range("AF5").Value = 500
$ws.Range("AG5").Value = 16
$ws.Range("AH5").Value = 500
$ws.Range("AI5").Value = 700
$ws.Range("AM5").Value = 700
$ws.Range("AN5").Value = 50

# Row 6
$ws.Range("I6").Value = 18
$ws.Range("N6").Value = 6
$ws.Range("P6").Value = 2.64
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("U6").Value = 1.75
$ws.Range("V6").Value = 1.05
$ws.Range("X6").Value = 1000
$ws.Range("AC6").Value = 17
$ws.Range("AF6").Value = 8.4

# Row 7
$ws.Range("I7").Value = 1.19
$ws.Range("J7").Value = 8.6
$ws.Range("O7").Value = 1.16
$ws.Range("R7").Value = 1.75
$ws.Range("S7").Value = 2.12
$ws.Range("T7").Value = 2.28
$ws.Range("V7").Value = 6.2
$ws.Range("Y7").Value = 21
$ws.Range("AA7").Value = 9.6
$ws.Range("AI7").Value = 160

# Row 8
$ws.Range("F8").Value = 2.92
$ws.Range("G8").Value = 3.1
$ws.Range("I8").Value = 3.3
$ws.Range("J8").Value = 2.7
$ws.Range("K8").Value = 2.86
$ws.Range("L8").Value = 1.68
$ws.Range("N8").Value = 2.32
$ws.Range("T8").Value = 2.3
$ws.Range("U8").Value = 1.68
$ws.Range("V8").Value = 1.43
$ws.Range("W8").Value = 1.47
$ws.Range("X8").Value = 6.8
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 9.4
$ws.Range("AC8").Value = 6.8
$ws.Range("AH8").Value = 32
$ws.Range("AJ8").Value = 1000
$ws.Range("AK8").Value = 1000

# Row 9
$ws.Range("F9").Value = 1.81
$ws.Range("G9").Value = 1.88
$ws.Range("H9").Value = 6
$ws.Range("I9").Value = 6.8
$ws.Range("J9").Value = 3.25
$ws.Range("K9").Value = 3.45
$ws.Range("L9").Value = 1.55
$ws.Range("O9").Value = 1.52
$ws.Range("P9").Value = 1.57
$ws.Range("Q9").Value = 2.5
$ws.Range("R9").Value = 1.21
$ws.Range("T9").Value = 2.28
$ws.Range("U9").Value = 1.7
$ws.Range("W9").Value = 2.12
$ws.Range("X9").Value = 10.5
$ws.Range("Y9").Value = 1000
$ws.Range("AD9").Value = 32
$ws.Range("AG9").Value = 23
$ws.Range("AH9").Value = 85
$ws.Range("AJ9").Value = 22
$ws.Range("AK9").Value = 75
$ws.Range("AN9").Value = 21

# Row 10
$ws.Range("F10").Value = 1.72
$ws.Range("G10").Value = 1.8
$ws.Range("H10").Value = 5.8
$ws.Range("K10").Value = 3.95
$ws.Range("L10").Value = 1.43
$ws.Range("N10").Value = 3.6
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = 1.33
$ws.Range("S10").Value = 3.5
$ws.Range("T10").Value = 1.9
$ws.Range("W10").Value = 2.24
$ws.Range("X10").Value = 21
$ws.Range("Z10").Value = 55
$ws.Range("AA10").Value = 700
$ws.Range("AB10").Value = 9.199999999999999
$ws.Range("AF10").Value = 22
$ws.Range("AI10").Value = 700
$ws.Range("AJ10").Value = 1000
$ws.Range("AL10").Value = 500
$ws.Range("AM10").Value = 500
$ws.Range("AN10").Value = 14

# Row 11
$ws.Range("I11").Value = 3.6
$ws.Range("L11").Value = 1.57
$ws.Range("N11").Value = 2.5
$ws.Range("O11").Value = 1.59
$ws.Range("P11").Value = 1.5
$ws.Range("R11").Value = 1.18
$ws.Range("U11").Value = 1.75
$ws.Range("V11").Value = 1.38
$ws.Range("Z11").Value = 26
$ws.Range("AA11").Value = 85
$ws.Range("AD11").Value = 18.5
$ws.Range("AH11").Value = 28
$ws.Range("AI11").Value = 95
$ws.Range("AK11").Value = 48
$ws.Range("AN11").Value = 55
$ws.Range("AO11").Value = 90

# Row 12
$ws.Range("I12").Value = 3.75
$ws.Range("J12").Value = 3.4
$ws.Range("K12").Value = 3.45
$ws.Range("N12").Value = 3.6
$ws.Range("S12").Value = 3.9
$ws.Range("T12").Value = 1.83
$ws.Range("U12").Value = 2.12
$ws.Range("V12").Value = 1.36
$ws.Range("Y12").Value = 13
$ws.Range("Z12").Value = 26
$ws.Range("AA12").Value = 70
$ws.Range("AB12").Value = 9.4
$ws.Range("AC12").Value = 7.4
$ws.Range("AE12").Value = 46
$ws.Range("AF12").Value = 14
$ws.Range("AG12").Value = 11
$ws.Range("AH12").Value = 18.5
$ws.Range("AJ12").Value = 28
$ws.Range("AK12").Value = 26
$ws.Range("AL12").Value = 44
$ws.Range("AM12").Value = 110
$ws.Range("AN12").Value = 20
$ws.Range("AO12").Value = 50

# Row 13
$ws.Range("F13").Value = 2.2
$ws.Range("G13").Value = 2.46
$ws.Range("H13").Value = 3.45
$ws.Range("I13").Value = 3.95
$ws.Range("J13").Value = 3.1
$ws.Range("K13").Value = 3.55
$ws.Range("S13").Value = 4.6
$ws.Range("U13").Value = 1.84
$ws.Range("V13").Value = 1.33
$ws.Range("W13").Value = 1.68
$ws.Range("X13").Value = 12.5
$ws.Range("Y13").Value = 13
$ws.Range("AC13").Value = 9
$ws.Range("AD13").Value = 19
$ws.Range("AE13").Value = 70
$ws.Range("AF13").Value = 16.5
$ws.Range("AI13").Value = 500

# Row 14
$ws.Range("F14").Value = 1.68
$ws.Range("N14").Value = 3.45
$ws.Range("Q14").Value = 2.1
$ws.Range("U14").Value = 1.87
$ws.Range("W14").Value = 2.4
$ws.Range("Y14").Value = 19.5
$ws.Range("AA14").Value = 190
$ws.Range("AD14").Value = 25
$ws.Range("AE14").Value = 110
$ws.Range("AH14").Value = 25
$ws.Range("AI14").Value = 110
$ws.Range("AM14").Value = 160
$ws.Range("AO14").Value = 150

# Row 15
$ws.Range("F15").Value = 1.47
$ws.Range("P15").Value = 1.84
$ws.Range("Q15").Value = 2.06
$ws.Range("R15").Value = 1.32
$ws.Range("S15").Value = 3.75
$ws.Range("Y15").Value = 29
$ws.Range("AA15").Value = 480
$ws.Range("AD15").Value = 42
$ws.Range("AL15").Value = 60
$ws.Range("AO15").Value = 400

# Row 16
$ws.Range("I16").Value = 10.5
$ws.Range("Y16").Value = 23
$ws.Range("AH16").Value = 38
